$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert "Example1Imp" right after "Example1Obs" (3rd sheet)
# ---------------------------------------------------------------------------
$obs1 = $wb.Worksheets.Item("Example1Obs")
$imp1 = $wb.Worksheets.Add($null, $obs1)
$imp1.Name = "Example1Imp"

# ---------------------------------------------------------------------------
# 2. Insert "Example2Imp" at the very end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$imp2 = $wb.Worksheets.Add($null, $lastSheet)
$imp2.Name = "Example2Imp"

# ---------------------------------------------------------------------------
# Helper: fill in the standard "Imp" sheet content (same data + labels on
# both Example1Imp / Example2Imp). Writing column A top-to-bottom first and
# then column B keeps new shared-string insertion order stable:
# TACSD, TACFrac, ESD, EFrac, SizeLimSD, SizeLimFrac, DLMtool generated,
# Perfect_Imp
# ---------------------------------------------------------------------------
function Fill-ImpSheet($ws) {
    $ws.Range("A1").Value = "Name"
    $ws.Range("A2").Value = "TACSD"
    $ws.Range("A3").Value = "TACFrac"
    $ws.Range("A4").Value = "ESD"
    $ws.Range("A5").Value = "EFrac"
    $ws.Range("A6").Value = "SizeLimSD"
    $ws.Range("A7").Value = "SizeLimFrac"
    $ws.Range("A8").Value = "Source"

    $ws.Range("B2").Value = 0
    $ws.Range("C2").Value = 0
    $ws.Range("B3").Value = 1
    $ws.Range("C3").Value = 1
    $ws.Range("B4").Value = 0
    $ws.Range("C4").Value = 0
    $ws.Range("B5").Value = 1
    $ws.Range("C5").Value = 1
    $ws.Range("B6").Value = 0
    $ws.Range("C6").Value = 0
    $ws.Range("B7").Value = 1
    $ws.Range("C7").Value = 1
    $ws.Range("B8").Value = "DLMtool generated"
    $ws.Range("B1").Value = "Perfect_Imp"

    # header / label formatting: Lucida Console 10pt black on solid white
    # fill, wrapped text - applied to the label column (A1:A8) and to the
    # "DLMtool generated" source cell (B8). Applied as two separate Range
    # operations (rather than a single multi-area "A1:A8,B8" union) so both
    # areas reliably pick up the formatting.
    $labels = $ws.Range("A1:A8")
    $labels.Font.Name = "Lucida Console"
    $labels.Font.Size = 10
    $labels.Font.Color = 0
    $labels.Interior.Color = 16777215
    $labels.WrapText = $true

    $source = $ws.Range("B8")
    $source.Font.Name = "Lucida Console"
    $source.Font.Size = 10
    $source.Font.Color = 0
    $source.Interior.Color = 16777215
    $source.WrapText = $true
}

Fill-ImpSheet $imp1
Fill-ImpSheet $imp2

# ---------------------------------------------------------------------------
# Row heights for the wrapped, longer labels on Example1Imp
# ---------------------------------------------------------------------------
$imp1.Rows.Item(6).RowHeight = 26.25
$imp1.Rows.Item(7).RowHeight = 26.25
$imp1.Rows.Item(8).RowHeight = 39

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$imp1.Columns.Item(2).ColumnWidth = 11.86

$imp2.Columns.Item(1).ColumnWidth = 13.71
$imp2.Columns.Item(2).ColumnWidth = 20.71

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
$imp2.Range("A1:D8").Select() | Out-Null
$imp1.Range("B12").Select() | Out-Null
$imp1.Activate() | Out-Null
